$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("N2").Formula = "=_xlfn.STDEV.S(C2:C8)"
$ws.Range("N3").Formula = "=_xlfn.STDEV.S(D2:D8)"
$ws.Range("N4").Formula = "=_xlfn.STDEV.S(E2:E8)"
$ws.Range("N5").Formula = "=_xlfn.STDEV.S(F2:F8)"

$ws.Range("N6").Select()
